$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark from its old location
#    (it used to sit right after the last character of "לוח תהיה
#    השתקפות" in the "Optional" bullet list).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Turn the "דוגמאות:" paragraph into "דוגמאות להשראה:" -- inserted
#    as its own run so the final layout is:
#       run1 "דוגמאות"
#       run2 " להשראה"
#       bookmark "_GoBack" (collapsed, right before the colon)
#       run3 ":"
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("דוגמאות:", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # $rng is now collapsed onto the found text "דוגמאות:"
    # Insert " להשראה" right before the trailing colon.
    $insertPoint = $d.Range($rng.End - 1, $rng.End - 1)
    $insertPoint.InsertBefore(" להשראה")

    # Force a run boundary between "דוגמאות" and " להשראה" (Word keeps
    # them as two separate runs even though the formatting is
    # identical) by toggling a character property on/off.
    $rng2 = $d.Content
    $found2 = $rng2.Find.Execute(" להשראה", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
    if ($found2) {
        $rng2.Bold = 1
        $rng2.Bold = 0
    }

    # Re-locate the full phrase and drop the (collapsed) "_GoBack"
    # bookmark immediately before the final colon -- this both marks
    # the last edit position and forces the trailing ":" into its own
    # run.
    $rng3 = $d.Content
    $found3 = $rng3.Find.Execute("דוגמאות להשראה:", $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "", 0)
    if ($found3) {
        $bmPoint = $d.Range($rng3.End - 1, $rng3.End - 1)
        $d.Bookmarks.Add("_GoBack", $bmPoint)
    }
}
